$wb = $excel.ActiveWorkbook
$wsBien = $wb.Worksheets.Item("Bien")
$wsServicio = $wb.Worksheets.Item("Servicio")

# =====================================================================
# Sheet "Bien" (sheet1): add a second "BIEN" row (row 4) and bump the
# quantity already present in row 3.
# =====================================================================

# Row 3: quantity value changes from 5 to 454
$wsBien.Range("E3").Value = 454

# Copy formatting from row 3 down to the new row 4 so fonts/fills/borders match
$wsBien.Range("A3:E3").Copy()
$wsBien.Range("A4:E4").PasteSpecial(-4122)

# Row 4 values
$wsBien.Range("A4").Value = "1 - BIEN"
$wsBien.Range("B4").Value = "BIEN 2"
$wsBien.Range("C4").Value = "DESCRIPCION BIEN 2"
$wsBien.Range("D4").Value = "5 - GRAMO"
$wsBien.Range("E4").Value = 5

# Quantity column now carries 2 decimal places
$wsBien.Range("E3:E4").NumberFormat = "0.00"

# Match row 3's custom row height on the newly added row 4
$wsBien.Rows.Item(4).RowHeight = $wsBien.Rows.Item(3).RowHeight

# Re-create the data validations so they cover both row 3 and row 4, and
# reorder them to land: list(D), list(A), decimal(E)
$wsBien.Range("D3:D4").Validation.Delete()
$wsBien.Range("D3:D4").Validation.Add(3, 1, 1, '"0 - NO APLICA,1 - METRO,2 - KILOMETRO,3 - CENTIMETRO,4 - KILOGRAMO,5 - GRAMO,6 - MESES,7 - DIAS,8 - AÑOS,9 - SEGUNDOS,10 - MINUTOS,11 - HORAS,12 - LITRO,13 - UNIDAD"', "0")
$wsBien.Range("D3:D4").Validation.ShowInput = $true
$wsBien.Range("D3:D4").Validation.ShowError = $true

$wsBien.Range("A3:A4").Validation.Delete()
$wsBien.Range("A3:A4").Validation.Add(3, 1, 1, '"1 - BIEN"', "0")
$wsBien.Range("A3:A4").Validation.IgnoreBlank = $false
$wsBien.Range("A3:A4").Validation.ShowInput = $true
$wsBien.Range("A3:A4").Validation.ShowError = $true

$wsBien.Range("E3:E4").Validation.Delete()
$wsBien.Range("E3:E4").Validation.Add(2, 1, 1, "0", "1.11111111111111E+015")
$wsBien.Range("E3:E4").Validation.ShowInput = $false
$wsBien.Range("E3:E4").Validation.ShowError = $true

# =====================================================================
# Sheet "Servicio" (sheet2): row 4 tweaks
# =====================================================================

# Unidad selection changes from "3 - CENTIMETRO" to "0 - NO APLICA"
$wsServicio.Range("D4").Value = "0 - NO APLICA"

# H4 (Cantidad del Servicio) becomes a decimal amount
$wsServicio.Range("H4").Value = 456465.45
$wsServicio.Range("H4").NumberFormat = "0.00"

# Re-apply the validations that were already present (E4, F4, G4, D4, A4)
# so their flags survive the round-trip, keeping their original relative
# order, then finally re-create H4's validation as "decimal" (moves to
# the end of the list, like in the target workbook).
$wsServicio.Range("E4").Validation.Delete()
$wsServicio.Range("E4").Validation.Add(1, 1, 1, "0", "1E+019")
$wsServicio.Range("E4").Validation.ShowInput = $true
$wsServicio.Range("E4").Validation.ShowError = $true

$wsServicio.Range("F4").Validation.Delete()
$wsServicio.Range("F4").Validation.Add(3, 1, 1, '"0,1,2,3,4,5,6,7,8,9,10,11,12"', "1E+019")
$wsServicio.Range("F4").Validation.ShowInput = $true
$wsServicio.Range("F4").Validation.ShowError = $true

$wsServicio.Range("G4").Validation.Delete()
$wsServicio.Range("G4").Validation.Add(3, 1, 1, '"0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30"', "1E+019")
$wsServicio.Range("G4").Validation.ShowInput = $true
$wsServicio.Range("G4").Validation.ShowError = $true

$wsServicio.Range("D4").Validation.Delete()
$wsServicio.Range("D4").Validation.Add(3, 1, 1, '"0 - NO APLICA,1 - METRO,2 - KILOMETRO,3 - CENTIMETRO,4 - KILOGRAMO,5 - GRAMO,6 - MESES,7 - DIAS,8 - AÑOS,9 - SEGUNDOS,10 - MINUTOS,11 - HORAS,12 - LITRO,13 - UNIDAD"', "0")
$wsServicio.Range("D4").Validation.ShowInput = $true
$wsServicio.Range("D4").Validation.ShowError = $true

$wsServicio.Range("A4").Validation.Delete()
$wsServicio.Range("A4").Validation.Add(3, 1, 1, '"2 - SERVICIO"', "0")
$wsServicio.Range("A4").Validation.IgnoreBlank = $false
$wsServicio.Range("A4").Validation.ShowInput = $true
$wsServicio.Range("A4").Validation.ShowError = $true

# H4 validation switches from whole numbers to decimal, and moves to the
# end of the dataValidations list
$wsServicio.Range("H4").Validation.Delete()
$wsServicio.Range("H4").Validation.Add(2, 1, 1, "0", "1.11111111111111E+015")
$wsServicio.Range("H4").Validation.ShowInput = $false
$wsServicio.Range("H4").Validation.ShowError = $true

Write-Host "edit applied"
